# "Maybe fixed typing bug" - updates to the lookups workbook:
#   * Companies: add a new active company "Test"
#   * Locations: BC is now also used by company "Test"
#   * AssetTypes: recolor the existing Cableway/BC/NHS row and add a new
#     Cableway/BC/Test row
#   * Settings: applyStatusColorsOnMap flips from FALSE to TRUE
$wb = $excel.ActiveWorkbook

# --- Companies sheet: add a new "Test" company (active) ---
$wsCompanies = $wb.Worksheets.Item("Companies")
$wsCompanies.Range("A3").Value = "Test"
# Copy the existing text "TRUE" cell (B2) so the new cell keeps the same
# text type instead of Excel auto-converting a literal TRUE/FALSE typed
# value into a boolean cell.
$wsCompanies.Range("B2").Copy()
$wsCompanies.Range("B3").PasteSpecial(-4163)

# --- Locations sheet: BC is now also used by company "Test" ---
$wsLocations = $wb.Worksheets.Item("Locations")
$wsLocations.Range("A3").Value = "BC"
$wsLocations.Range("B3").Value = "Test"

# --- AssetTypes sheet: update existing color, add new row for Test company ---
$wsAssetTypes = $wb.Worksheets.Item("AssetTypes")
$wsAssetTypes.Range("D2").Value = "#cd197e"
$wsAssetTypes.Range("A3").Value = "Cableway"
$wsAssetTypes.Range("B3").Value = "BC"
$wsAssetTypes.Range("C3").Value = "Test"
$wsAssetTypes.Range("D3").Value = "#196ce1"

# --- Settings sheet: flip applyStatusColorsOnMap from FALSE to TRUE ---
$wsSettings = $wb.Worksheets.Item("Settings")
$wsCompanies.Range("B2").Copy()
$wsSettings.Range("B2").PasteSpecial(-4163)
